# "minor fixes and corrections in the db"
#
# - add a new "Turkey" row (row 32) to the demand_storage_kWh sheet,
#   shifting the existing "United Kingdom" row down to row 33
#   (formulas that reference the UK row auto-adjust to the new row number)
# - fix Serbia's (row 26) Total-2021 figure, which was recorded as 0
# - flag both Turkey and Serbia with a "*" footnote marker in column H

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("demand_storage_kWh")

# Insert a new row above the old row 32 (United Kingdom). Excel shifts
# United Kingdom down to row 33 and rewrites any formulas that pointed at
# row 32 (e.g. Ireland's E17 "=0.7*E32" becomes "=0.7*E33") automatically.
$ws.Rows("32:32").Insert()

# New row 32: Turkey
$ws.Range("A32").Value = "Turkey"
$ws.Range("B32").Value = 344541000000
$ws.Range("E32").Value = 0.45
$ws.Range("H32").Value = "*"

# Serbia (row 26): correct the Total 2021 value and add the footnote marker
$ws.Range("B26").Value = 12867000000
$ws.Range("H26").Value = "*"

# Leave the selection where the author left it
$ws.Range("E33").Select()
